# Adds a new task "Validate Token" to the task list on Sheet1 (row 33,
# following the existing every-other-row layout), and moves the view /
# selection to reflect where the author was working afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task row (B33), mirrors the sparse row pattern already used for
# the rest of the list (last existing entry was B31).
$ws.Range("B33").Value = "Validate Token"

# Reflect the new selection/scroll position from the saved view.
[void]$ws.Range("A30").Select()
[void]$ws.Range("E33").Select()
